$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Color helper values (Excel Font.Color expects 0xBBGGRR packed integer)
$blue1 = 0xBC*65536 + 0x75*256 + 0x1B   # FF1B75BC
$blue2 = 0xB3*65536 + 0x66*256 + 0x00   # FF0066B3

# "Separate Weapon Animation from Hero Animation" row is now DONE; add D5 = DONE
# and recolor A5 text
$ws.Range("D5").Value = "DONE"
$ws.Range("A5").Font.Color = $blue1

# Recolor the "Level up/level down/reset sprite items..." row (row 9) and its DONE marker
$ws.Range("A9").Font.Color = $blue2
$ws.Range("D9").Font.Color = $blue2

# Add new backlog rows
$ws.Range("A11").Value = "Scripts for items"
$ws.Range("A12").Value = "Find better way to show shield in DodgeCombatAnimation"

$ws.Range("A13").Value = "Add jank AI"
$ws.Range("D13").Value = "DONE"
$ws.Range("A13").Font.Color = $blue1

$ws.Range("A14").Value = "Add wander approach type"
$ws.Range("D14").Value = "DONE"
$ws.Range("A14").Font.Color = $blue1

$ws.Range("A15").Value = "Add AI vision"
$ws.Range("D15").Value = "DONE"
$ws.Range("A15").Font.Color = $blue1

# Update active selection to D6, matching the saved view state
[void]$ws.Range("D6").Select()
